$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(39)
$srcRange = $src.Range("A1:D51")
$srcRange.Copy()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-07-31"
$dst = $ws.Range("A1:D51")
$dst.PasteSpecial(-4122)

$data = New-Object 'object[,]' 51,4
$data[0,0] = 'rank'
$data[0,1] = 'title'
$data[0,2] = 'author'
$data[0,3] = 'latest_episode'
$data[1,0] = 1
$data[1,1] = '生徒会にも穴はある！'
$data[1,2] = 'むちまろ'
$data[1,3] = '第132話	ありす大ピンチ！（バッドエンド編）'
$data[2,0] = 2
$data[2,1] = '実は俺、最強でした？'
$data[2,2] = '原作：澄守 彩 漫画：高橋 愛'
$data[2,3] = 'おまけ64'
$data[3,0] = 3
$data[3,1] = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$data[3,2] = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$data[3,3] = '第32話 独身貴族は森で写真を撮る（3）'
$data[4,0] = 4
$data[4,1] = '異世界食堂　洋食のねこや'
$data[4,2] = '犬塚惇平(ヒーロー文庫／イマジカインフォス)(原作) ヤミザワ(漫画) モロザワ(漫画) エナミカツミ(キャラクター原案)'
$data[4,3] = '第40話②'
$data[5,0] = 5
$data[5,1] = 'オタクに優しいギャルはいない!?'
$data[5,2] = 'のりしろちゃん 魚住さかな'
$data[5,3] = '【#150】ひび割れ'
$data[6,0] = 6
$data[6,1] = '生徒会役員共'
$data[6,2] = '氏家ト全'
$data[6,3] = '#407'
$data[7,0] = 7
$data[7,1] = '脇役に転生したはずが、いつの間にか伝説の錬金術師になってた～仲間たちが英雄でも俺は支援職なんだが～'
$data[7,2] = '神無月みり 相野 仁'
$data[7,3] = '第２６話　脇役、目論見を見破る（３）'
$data[8,0] = 8
$data[8,1] = '♀ガキとおじさん'
$data[8,2] = 'サラマンダ(著者)'
$data[8,3] = '第29話'
$data[9,0] = 9
$data[9,1] = 'カナン様はあくまでチョロい'
$data[9,2] = 'nonco'
$data[9,3] = '第144話	リリイの仕返し恩返し'
$data[10,0] = 10
$data[10,1] = 'ブレイド＆バスタード'
$data[10,2] = '漫画/楓月 誠 原作/蝸牛くも キャラクター原案/so-bin'
$data[10,3] = '第10話（2）'
$data[11,0] = 11
$data[11,1] = 'すべての人類を破壊する。それらは再生できない。'
$data[11,2] = '横田卓馬(漫画) 伊瀬勝良(原作)'
$data[11,3] = '第67話その1'
$data[12,0] = 12
$data[12,1] = '色欲無双 ～変態スキルが暴走してヤリサーから追放された俺は、はからずも淫靡な力で最強になる～'
$data[12,2] = 'あいのひとし 笠原巴 三九呂'
$data[12,3] = '第1話 ヤリサーを追放！？'
$data[13,0] = 13
$data[13,1] = 'よわよわ先生'
$data[13,2] = '福地カミオ'
$data[13,3] = '第109話	ほどほどのバレンタインデー'
$data[14,0] = 14
$data[14,1] = '世界最速のレベルアップ'
$data[14,2] = '鈴見敦(漫画) 八又ナガト(原作) fame(キャラクター原案)'
$data[14,3] = '第48話②'
$data[15,0] = 15
$data[15,1] = '異世界居酒屋「のぶ」'
$data[15,2] = '蝉川夏哉(原作) ヴァージニア二等兵(漫画) 転(キャラクター原案)'
$data[15,3] = '第123話'
$data[16,0] = 16
$data[16,1] = 'おかけになった呪文は現在使われておりません'
$data[16,2] = 'ロケット商会 天宮ケイリ'
$data[16,3] = '第1話　婚活の呪文'
$data[17,0] = 17
$data[17,1] = '善人おっさん、生まれ変わったらSSSランク人生が確定した'
$data[17,2] = '原作／三木なずな 漫画／ゆづましろ キャラクター原案／伍長'
$data[17,3] = '祝！単行本11巻発売！特別イラスト'
$data[18,0] = 18
$data[18,1] = 'ギャルゲーマーに褒められたい'
$data[18,2] = 'げしゅまろ(著者)'
$data[18,3] = '45話'
$data[19,0] = 19
$data[19,1] = 'やり直し令嬢は竜帝陛下を攻略中'
$data[19,2] = '柚アンコ(漫画) 永瀬さらさ（角川ビーンズ文庫）(原作) 藤未都也(キャラクター原案)'
$data[19,3] = 'Episode40.5'
$data[20,0] = 20
$data[20,1] = '百瀬アキラの初恋破綻中。'
$data[20,2] = '晴川シンタ'
$data[20,3] = '第33話 あくまで公務を執行中。'
$data[21,0] = 21
$data[21,1] = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$data[21,2] = '漫画/すたひろ 原作/Y.A'
$data[21,3] = 'chapter66【35話①】'
$data[22,0] = 22
$data[22,1] = 'え、社内システム全てワンオペしている私を解雇ですか？'
$data[22,2] = '漫画：伊於 原作：下城米雪 キャラクター原案：icchi'
$data[22,3] = '4巻発売告知漫画'
$data[23,0] = 23
$data[23,1] = 'ガリ勉くんと裏アカさん　散々お世話になっているエロ系裏垢女子の正体がクラスのアイドルだった件'
$data[23,2] = '花咲まひる(著者) 鈴木えんぺら(原作) 小花雪(キャラクター原案)'
$data[23,3] = '第9話①'
$data[24,0] = 24
$data[24,1] = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$data[24,2] = '戸賀 環 坂木持丸 riritto'
$data[24,3] = '第50話①　祝われた家を探索してみた'
$data[25,0] = 25
$data[25,1] = '女子高生の無駄づかい'
$data[25,2] = 'ビーノ(著者)'
$data[25,3] = '第133話　てこいれ'
$data[26,0] = 26
$data[26,1] = '奈落の底で生活して早三年、当時『白魔道士』だった私は『聖魔女』になっていた'
$data[26,2] = '雪葵ムラサキ(漫画) tani(原作) れんた(キャラクター原案)'
$data[26,3] = '第7話①'
$data[27,0] = 27
$data[27,1] = '妹はカノジョにできないのに'
$data[27,2] = 'ちくわ。(作画) 鏡遊(原作) 三九呂(キャラクターデザイン)'
$data[27,3] = '第29話②'
$data[28,0] = 28
$data[28,1] = '黒月のイェルクナハト'
$data[28,2] = 'スズモトコウ'
$data[28,3] = '第6話	生きる意味'
$data[29,0] = 29
$data[29,1] = 'ひとりぼっちの異世界攻略'
$data[29,2] = 'びび（漫画） 五示正司（原作）'
$data[29,3] = '第229話　最初は良い感じだったよ…？'
$data[30,0] = 30
$data[30,1] = '絶対死なないステラ姫'
$data[30,2] = '光永康則 大高稲'
$data[30,3] = '第１４話　絶対旅立たない（３）'
$data[31,0] = 31
$data[31,1] = '帰ってください！ 阿久津さん'
$data[31,2] = '長岡太一(著者)'
$data[31,3] = '第193話'
$data[32,0] = 32
$data[32,1] = 'ぽんドロイド！ はまさん'
$data[32,2] = 'はれやまはれぞう(著者)'
$data[32,3] = '第5話'
$data[33,0] = 33
$data[33,1] = 'アンゴルモア 元寇合戦記　【博多編】'
$data[33,2] = 'たかぎ七彦(著者)'
$data[33,3] = '第四十四話その七'
$data[34,0] = 34
$data[34,1] = '「おかえり、パパ」'
$data[34,2] = '蝉丸'
$data[34,3] = '第26話　家族'
$data[35,0] = 35
$data[35,1] = '最強の少年聖騎士、転生者を狩る'
$data[35,2] = '作画：御塩 原作：宇奈木ユラ'
$data[35,3] = '第6話(2)'
$data[36,0] = 36
$data[36,1] = 'ロードマギアの弟子'
$data[36,2] = 'FLIPFLOPs'
$data[36,3] = '第19話 魔術師の戦い (後編)'
$data[37,0] = 37
$data[37,1] = '姫ヶ崎櫻子は今日も不憫可愛い'
$data[37,2] = '安田剛助(著者)'
$data[37,3] = '第49話'
$data[38,0] = 38
$data[38,1] = 'ポーション、わが身を助ける'
$data[38,2] = '戸部淑 岩船晶'
$data[38,3] = '第1話'
$data[39,0] = 39
$data[39,1] = '異世界おじさん'
$data[39,2] = '殆ど死んでいる(著者)'
$data[39,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[40,0] = 40
$data[40,1] = '魔術師クノンは見えている'
$data[40,2] = 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)'
$data[40,3] = '第39話①'
$data[41,0] = 41
$data[41,1] = 'カッコウの許嫁'
$data[41,2] = '吉河美希'
$data[41,3] = '第257話	また将来の選択肢が増えちゃった!'
$data[42,0] = 42
$data[42,1] = '色憑くモノクローム'
$data[42,2] = '内山敦司'
$data[42,3] = '第39話	湧き上がる熱情'
$data[43,0] = 43
$data[43,1] = '江戸前エルフ'
$data[43,2] = '樋口彰彦'
$data[43,3] = '#116'
$data[44,0] = 44
$data[44,1] = 'ダウナー系お姉さんに毎日カスの嘘を流し込まれる話'
$data[44,2] = '生倉のゑる(著者) はるばーど屋(原作者)'
$data[44,3] = '11話 おまけ'
$data[45,0] = 45
$data[45,1] = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$data[45,2] = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$data[45,3] = '第５０話　雌雄を決する器用貧乏（４）'
$data[46,0] = 46
$data[46,1] = '魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～'
$data[46,2] = '漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ'
$data[46,3] = '第47話 魔導具師とつながれたもの③'
$data[47,0] = 47
$data[47,1] = '最後のエルフ'
$data[47,2] = 'サワノアキラ（漫画）'
$data[47,3] = '第9章　竜の眠る地（後編）'
$data[48,0] = 48
$data[48,1] = '帰ってきた元勇者'
$data[48,2] = '漫画：なるさわ景 原作：ニシ キャラクター原案：米白粕'
$data[48,3] = '第31話(3)'
$data[49,0] = 49
$data[49,1] = 'ゆめねこねくと'
$data[49,2] = '澤田コウ'
$data[49,3] = '第38こねくと	ゆめねこねくと'
$data[50,0] = 50
$data[50,1] = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$data[50,2] = '作画：マエD 原作：新人'
$data[50,3] = '第5話(2)'

$ws.Range("A1:D51").Value = $data
